$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 2/3 text tweaks first, so the shared-string table rebuilds in the right order ---
# Renaming A3 drops the old "sales_order" string and re-adds "sales_orders" right after the
# still-used strings (FeatureFileName, ResellerPO, ResellerBCN, aged_orders).
$ws.Range("A3").Value = "sales_orders"

# --- 2. Clear the old sequence-number cells (C2, C3) ---
$ws.Range("C2").ClearContents() | Out-Null
$ws.Range("C3").ClearContents() | Out-Null

# --- 3. Capture the formatting of the existing header cells before we overwrite them ---
# A1/B1 use the "full border" bold header style; C1 uses the "left/right border only" bold header style.
# We copy that format onto the new header range so the resulting style table reuses the same
# font/border combination as the original file.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("C1:G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- 4. Rewrite the new header/data cells, interleaved in this exact order so the shared-string
# table is rebuilt with the same ordering as the target file ---
$ws.Range("C1").Value = "IMOrderNo"
$ws.Range("C3").Value = "20-RP39N-11"
$ws.Range("D1").Value = "OrderType"
$ws.Range("D3").Value = "Stock"
$ws.Range("E1").Value = "ResellerPO"
$ws.Range("F1").Value = "VendorName"
$ws.Range("G1").Value = "OrderStatus"
$ws.Range("E3").Value = "PO12345"
$ws.Range("F3").Value = "MICROSOFT HARDWARE"
$ws.Range("G3").Value = "Order Hold"

# --- 5. Column widths for the new/changed columns ---
# (target widths, in "characters": C=13, D/E=14.43, F=25.14, G=12.14)
$ws.Columns.Item(3).ColumnWidth = 12.1875
$ws.Columns.Item(4).ColumnWidth = 13.59375
$ws.Columns.Item(5).ColumnWidth = 13.59375
$ws.Columns.Item(6).ColumnWidth = 24.375
$ws.Columns.Item(7).ColumnWidth = 11.25

# --- 7. Selection matches the new active cell position ---
$ws.Range("F11").Select() | Out-Null
